$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.458.37"
$ws.Range("E2").Value = "  +2.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.013.68"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.59"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("E6").Value = "  +2.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.71"
$ws.Range("E7").Value = "  +5.18%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0807"
$ws.Range("E10").Value = "  +0.53%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("E12").Value = "  +8.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.95"
$ws.Range("E13").Value = "  +3.95%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.855"
$ws.Range("E14").Value = "  +2.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.306.25"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.009.17"
$ws.Range("E17").Value = "  +2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.406.57"
$ws.Range("E18").Value = "  +2.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.65"
$ws.Range("E19").Value = "  +1.16%  "
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("E21").Value = "  +3.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.50"
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  +3.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.41"
$ws.Range("E27").Value = "  +1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.07"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.77"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("E30").Value = "  +18.99%  "
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.88"
$ws.Range("E32").Value = "  +3.67%  "
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("E35").Value = "  +5.54%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.82"
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.59"
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0984"
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.94"
$ws.Range("E41").Value = "  +1.31%  "
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.88"
$ws.Range("E44").Value = "  +5.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.73"
$ws.Range("E45").Value = "  +4.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.384.26"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.31"
$ws.Range("E48").Value = "  +2.14%  "
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.93"
$ws.Range("E50").Value = "  +7.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.99"
$ws.Range("E51").Value = "  +12.01%  "
